$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.172.12"
$ws.Range("D3").Value = "1.785.41"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.95"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.548"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.07"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.293"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0688"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0949"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D12").Value = "2.042.35"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.822.97"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.98"
$ws.Range("E14").Value = "  -3.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.625"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").Value = "34.144.16"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.19"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.72"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.66"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "0.0₃0796"
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.00"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.14"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.19"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0519"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.22"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.75"
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.73"
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").Value = "1.444.85"
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.63"
$ws.Range("E36").Value = "  +11.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.656"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.06"
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.07"
$ws.Range("E41").Value = "  +5.29%  "
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.916"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.09"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("D48").Value = "1.942.40"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.93"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("E50").Value = "  -5.36%  "
$ws.Range("E51").Value = "  +0.22%  "
